# Add MinimalInformationGain to stop logic and make stop logic usable
#
# - Options!B3 ("Minimal information gain") changes from 10 to 0.05
# - Selections are updated to reflect where the user left the cursor:
#     Survey  sheet: F8  -> D8
#     Options sheet: A4  -> A3:B3 (and Options stays the active/selected tab)

$wb = $excel.ActiveWorkbook

$wsSurvey  = $wb.Worksheets.Item("Survey")
$wsOptions = $wb.Worksheets.Item("Options")

# Update the value that drives the stop logic's minimal information gain.
$wsOptions.Range("B3").Value = 0.05

# Leave the cursor on the Survey sheet where the author last clicked (D8),
# without leaving that sheet active.
$wsSurvey.Range("D8").Select() | Out-Null

# Options is the active/visible sheet; select A3:B3 with A3 as the active cell.
$wsOptions.Activate() | Out-Null
$wsOptions.Range("A3:B3").Select() | Out-Null
